$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 70, shifting old row 70 down to row 71.
$ws.Rows.Item(70).Insert()

# New row 70 becomes a copy of the (about to be overwritten) row 69 data
# (the original row 69 values, before the update below is applied).
$ws.Range("A70:R70").Value = $ws.Range("A69:R69").Value()

# Copy the date cell format (style) used on column D down to the new D70 cell.
$ws.Range("D69").Copy()
$ws.Range("D70").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Now update row 69 with its new values per the diff.
$ws.Range("D69").Value = 45021
$ws.Range("J69").Value = 700
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 7000
$ws.Range("M69").Value = 6500
$ws.Range("P69").Value = 108
